$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @{
    2 = @(0.9040748851485086, 1557358.170258267, 1247.941573255041, 822.0702927179702, 0.437571462660277)
    3 = @(0.9599677548029431, 649929.3144234334, 806.1819363043514, 400.8446421950315, 0.1130639482868977)
    4 = @(0.9524841413354613, 771426.8658694846, 878.3090947209215, 435.6455320726733, 0.1094570516174347)
    5 = @(0.9821133439625952, 290392.4583418807, 538.8807459372442, 274.6806620365182, 0.07076930797172286)
    6 = @(0.9821731157673121, 289422.0544119907, 537.9796040854994, 276.939906137627, 0.07359096544436548)
    7 = @(0.9759246594843534, 390866.6495926457, 625.1932897853636, 322.883120725276, 0.08876225634328097)
    8 = @(0.9658102259828125, 555075.946350279, 745.0341913967969, 513.6313806061189, 0.2503237323057665)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 3).Value = $vals[0]
    $ws.Cells.Item($row, 4).Value = $vals[1]
    $ws.Cells.Item($row, 5).Value = $vals[2]
    $ws.Cells.Item($row, 6).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
